$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Row 3
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

# Row 4
$ws.Range("D4").Value = 0.9993305975431873
$ws.Range("E4").Value = 0.9993305975431873

# Row 5
$ws.Range("D5").Value = [double]"2.518342402627153E-06"
$ws.Range("E5").Value = [double]"2.518342402627153E-06"

# Row 6
$ws.Range("D6").Value = 0.004933958571001721
$ws.Range("E6").Value = 0.004933958571001721

# Row 7
$ws.Range("D7").Value = 0.9999999999999705
$ws.Range("E7").Value = [double]"2.953193245502916E-14"

# Row 8
$ws.Range("D8").Value = 0.0002026292759377265
$ws.Range("E8").Value = 0.9997973707240623

# Row 9
$ws.Range("D9").Value = 0.9999999988089914
$ws.Range("E9").Value = [double]"1.191008625056611E-09"

# Row 10
$ws.Range("D10").Value = [double]"2.086481389027562E-18"

# Row 11
$ws.Range("D11").Value = [double]"4.154109954665919E-45"
$ws.Range("F11").Value = 26.42608642578125
